$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "468032"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = ""
